$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize connector words ("de", "del", "el", "la", "los") within place names
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("A27").Value = "Ciudad De México"
$ws.Range("A36").Value = "Estado De México"
$ws.Range("B36").Value = "Ecatepec De Morelos"
$ws.Range("B38").Value = "San Felipe Del Progreso"
$ws.Range("B49").Value = "San Francisco Del Rincón"
$ws.Range("B51").Value = "Acapulco De Juárez"
$ws.Range("B54").Value = "Atenango Del Río"
$ws.Range("B57").Value = "Huitzuco De Los Figueroa"
$ws.Range("B58").Value = "Iguala De La Independencia"
$ws.Range("B59").Value = "Zihuatanejo De Azueta"
$ws.Range("B64").Value = "Taxco De Alarcón"
$ws.Range("B65").Value = "Técpan De Galeana"
$ws.Range("B68").Value = "Cuautepec De Hinojosa"
$ws.Range("B69").Value = "Huasca De Ocampo"
$ws.Range("B71").Value = "Tezontepec De Aldama"
$ws.Range("B72").Value = "Tulancingo De Bravo"
$ws.Range("B75").Value = "Autlán De Navarro"
$ws.Range("B81").Value = "Lagos De Moreno"
$ws.Range("B84").Value = "San Diego De Alejandría"
$ws.Range("B85").Value = "San Miguel El Alto"
$ws.Range("B87").Value = "Tlajomulco De Zúñiga"
$ws.Range("B113").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B115").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B124").Value = "Cuetzalan Del Progreso"
$ws.Range("B125").Value = "Huitzilan De Serdán"
$ws.Range("B128").Value = "Tetela De Ocampo"
$ws.Range("B133").Value = "San Juan Del Río"
$ws.Range("B158").Value = "Cosamaloapan De Carpio"
$ws.Range("B159").Value = "Cosautlán De Carvajal"
$ws.Range("B164").Value = "Soledad De Doblado"
$ws.Range("B167").Value = "Vega De Alatorre"
$ws.Range("B175").Value = "Villa De Cos"

# Correct floating point value
$ws.Range("D89").Value = 0.09090909090909093

# Remove trailing footer rows (previously 180-184), leaving a blank row 179
$ws.Rows("180:184").Delete()
